# Add new header columns E:I ("Gender", "Birth-day", "Birth-month",
# "Birth-year", "Height") to the right of the existing "Password" column,
# and apply bold / centered / thin-boxed header formatting across the
# whole header row A1:I1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new header cells -------------------------------------------------
$ws.Range("E1").Value = "Gender"
$ws.Range("F1").Value = "Birth-day"
$ws.Range("G1").Value = "Birth-month"
$ws.Range("H1").Value = "Birth-year"
$ws.Range("I1").Value = "Height"

# --- header formatting --------------------------------------------------
# Build the target look (bold font, centered/top aligned, thin box border)
# on a scratch cell far away from the used range, then copy *only* the
# formatting onto the whole header row in a single paste so every cell in
# A1:I1 converges on the very same style record instead of accumulating a
# new one per cell/property.
$scratch = $ws.Range("Z100")
$scratch.Font.Bold = $true
$scratch.HorizontalAlignment = -4108   # xlCenter
$scratch.VerticalAlignment = -4160     # xlTop
$scratch.Borders.LineStyle = 1         # xlContinuous

$scratch.Copy()
$ws.Range("A1:I1").PasteSpecial(-4122) # xlPasteFormats
$scratch.Clear()

$null = $ws.Range("A1:I1").Select()
